$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.382.54"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.822.40"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5136"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3930"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07658"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.107"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.99"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.264"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.495"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "1.825.30"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.32"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001094"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06675"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.125"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "28.395.88"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.255"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.77"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "2.033.43"
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.108"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1089"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.648"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.655"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07110"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2207"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02324"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.153"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.781"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6243"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.172"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.391"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5871"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.62"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  -0.25%  "
